$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp text in A1
$ws.Range("A1").Value = "Datos actualizados a 23 de Abril de 2020 a las 07:22"

# Rusia (row 13): update F13
$ws.Range("F13").Value = 2300

# India (row 20): update B20, C20, E20
$ws.Range("B20").Value = 21797
$ws.Range("C20").Value = 427
$ws.Range("E20").Value = 16740

# Pakistan (row 32): update F32
$ws.Range("F32").Value = 60

# Lituania (row 75): update B75, C75, D75, E75
$ws.Range("B75").Value = 1398
$ws.Range("C75").Value = 28
$ws.Range("D75").Value = 399
$ws.Range("E75").Value = 961

# Bulgaria overtakes Hong Kong in the ranking, so the two rows swap places.
# Row 84 becomes Bulgaria (with its refreshed totals) and row 85 becomes
# Hong Kong (keeping the totals that used to sit in row 84).
$ws.Range("A84").Value = "Bulgaria"
$ws.Range("B84").Value = 1081
$ws.Range("C84").Value = 57
$ws.Range("D84").Value = 190
$ws.Range("E84").Value = 841
$ws.Range("F84").Value = 37
$ws.Range("G84").Value = 1
$ws.Range("H84").Value = 50

$ws.Range("A85").Value = "Hong Kong"
$ws.Range("B85").Value = 1034
$ws.Range("C85").Value = 0
$ws.Range("D85").Value = 678
$ws.Range("E85").Value = 352
$ws.Range("F85").Value = 8
$ws.Range("G85").Value = 0
$ws.Range("H85").Value = 4
